$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new nodes (previously mispositioned further out, now corrected to F1/G1)
$ws.Range("F1").Value = "Manguera3"
$ws.Range("G1").Value = "Elemento4"

# Copy formatting from the existing "manguera" style cell (D1) and "elemento" style cell (E1)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Update selection to reflect corrected active cell
$ws.Range("I5").Select() | Out-Null
